$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44488
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("Q3").Value = "$/bandeja 5 kilos"
$ws.Range("R3").Value = "La Ligua"
$ws.Range("S3").Value = 2400
$ws.Range("T3").Value = 5

# Row 4
$ws.Range("D4").Value = 44483
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 35
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 10000
$ws.Range("Q4").Value = "$/bandeja 5 kilos"
$ws.Range("R4").Value = "Provincia de Quillota"
$ws.Range("S4").Value = 2000
$ws.Range("T4").Value = 5

# Row 5
$ws.Range("D5").Value = 44515
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 28000
$ws.Range("O5").Value = 28000
$ws.Range("P5").Value = 28000
$ws.Range("Q5").Value = "$/bandeja 10 kilos"
$ws.Range("R5").Value = "Provincia de Los Andes"
$ws.Range("S5").Value = 2800
$ws.Range("T5").Value = 10

# Row 6
$ws.Range("D6").Value = 44519
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 28000
$ws.Range("O6").Value = 28000
$ws.Range("P6").Value = 28000
$ws.Range("Q6").Value = "$/bandeja 10 kilos"
$ws.Range("R6").Value = "Provincia de Quillota"
$ws.Range("S6").Value = 2800
$ws.Range("T6").Value = 10

# Row 7
$ws.Range("D7").Value = 44496
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 55
$ws.Range("N7").Value = 28000
$ws.Range("O7").Value = 28000
$ws.Range("P7").Value = 28000
$ws.Range("Q7").Value = "$/bandeja 10 kilos"
$ws.Range("R7").Value = "Provincia de Quillota"
$ws.Range("S7").Value = 2800
$ws.Range("T7").Value = 10

# Row 8
$ws.Range("D8").Value = 44511
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 45
$ws.Range("N8").Value = 28000
$ws.Range("O8").Value = 28000
$ws.Range("P8").Value = 28000
$ws.Range("Q8").Value = "$/bandeja 10 kilos"
$ws.Range("R8").Value = "Provincia de Los Andes"
$ws.Range("S8").Value = 2800
$ws.Range("T8").Value = 10

# Row 9
$ws.Range("D9").Value = 44511
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 45
$ws.Range("N9").Value = 3200
$ws.Range("O9").Value = 3200
$ws.Range("P9").Value = 3200
$ws.Range("Q9").Value = "$/bandeja 10 kilos"
$ws.Range("R9").Value = "Provincia de Quillota"
$ws.Range("S9").Value = 320
$ws.Range("T9").Value = 10

# Row 10
$ws.Range("D10").Value = 44466
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 11000
$ws.Range("O10").Value = 11000
$ws.Range("P10").Value = 11000
$ws.Range("Q10").Value = "$/bandeja 5 kilos"
$ws.Range("R10").Value = "La Ligua"
$ws.Range("S10").Value = 2200
$ws.Range("T10").Value = 5

# Row 11
$ws.Range("D11").Value = 44166
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 20
$ws.Range("N11").Value = 12000
$ws.Range("O11").Value = 12000
$ws.Range("P11").Value = 12000
$ws.Range("Q11").Value = "$/caja 18 kilos"
$ws.Range("R11").Value = "La Ligua"
$ws.Range("S11").Value = 667
$ws.Range("T11").Value = 18

Write-Output "Applied weekly fruit/vegetable update"
